$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("1:1").Insert(-4121)
$ws.Rows("62:62").Cut($ws.Rows("1:1"))
